# Final commit of upload excel file
# - Fill in the missing FirstName for row 2 ("rohan")
# - Update email + hobbies text for row 2
# - Update FirstName, Street and hobbies text for row 3
# - Bump header/data row heights slightly
# - Normalize the Pincode/Phone column font color to explicit black

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits -------------------------------------------------
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Row height bump (18.75 -> 19.5) on the header + both data rows ---
$ws.Rows("1:3").RowHeight = 19.5

# --- Explicit black font color for the Pincode/Phone numeric columns --
$ws.Range("I2").Font.Color = 0
$ws.Range("I3").Font.Color = 0
$ws.Range("K2").Font.Color = 0
$ws.Range("K3").Font.Color = 0
